# Auto-update draw results: append the 2025-11-07 Pick 3 draw as a new row
# (row 52) at the bottom of the "Results" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (date-like text "2025-11-07") and C (plain digit text "251107")
# would otherwise be auto-coerced by Excel into a date serial / number.
# Pre-format them as Text so the values are stored as literal strings,
# matching every other row in this column.
$ws.Range("A52:A52").NumberFormat = "@"
$ws.Range("C52:C52").NumberFormat = "@"

$ws.Range("A52").Value = "2025-11-07"
$ws.Range("B52").Value = "Pick 3"
$ws.Range("C52").Value = "251107"
$ws.Range("D52").Value = "8-2-1"
$ws.Range("E52").Value = "2025-11-07T21:38:43.753+04:00"
